# Update gh-pages to output generated at 456a3b4
#
# Applies the same set of edits to both the "展览" sheet (which had 13 data
# rows, 2..13) and the "全部类型" sheet (which had 14 data rows, 2..14,
# because it also contains a "演出" row that "展览" doesn't have). Both
# sheets get:
#   - G2: 100 -> "不可售" (now text, not numeric)
#   - F5: 16 -> 18
#   - F6: 3197 -> 3212
#   - F7: 2096 -> 2098
#   - A brand-new event ("南宁·首届童话梦境Lolita茶会") inserted right
#     before "南宁·AB动漫游戏嘉年华", pushing every later row down by one
#   - "南宁·AB动漫游戏嘉年华"'s 想去人数 (F column): 1188 -> 1190
#   - "南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）"'s 想去人数: 1102 -> 1123
#   - "南宁·蔚蓝档案only"'s 想去人数: 90 -> 91
#   - column A renumbered sequentially for the extra row

$wb = $excel.ActiveWorkbook

function Update-ExpoSheet {
    param(
        [object]$ws,
        [int]$insertAt   # 1-based row index where the new "Lolita" row should end up
    )

    # --- simple value tweaks that do not move any rows -------------------
    $ws.Cells.Item(2, 7).Value = "不可售"

    $ws.Cells.Item(5, 6).Value = 18
    $ws.Cells.Item(6, 6).Value = 3212
    $ws.Cells.Item(7, 6).Value = 2098

    # --- insert the new row, shifting everything below it down by one ----
    $ws.Rows.Item($insertAt).Insert()

    $newRowNumber = $insertAt - 1   # sequence number shown in column A

    # Insert() doesn't reliably carry the bordered/bold/centered style that
    # every other column-A cell uses, so copy it explicitly from the row
    # just above (still untouched - same style as every other data row).
    $ws.Cells.Item($insertAt - 1, 1).Copy()
    $ws.Cells.Item($insertAt, 1).PasteSpecial(-4122)

    $ws.Cells.Item($insertAt, 1).Value = $newRowNumber
    $ws.Cells.Item($insertAt, 2).Value = "'2024-07-06"
    $ws.Cells.Item($insertAt, 3).Value = "南宁·首届童话梦境Lolita茶会"
    $ws.Cells.Item($insertAt, 4).Value = "明秀东路157号 利泰国际大酒店"
    $ws.Cells.Item($insertAt, 5).Value = "2024.07.06 13:00-07.06 17:00"
    $ws.Cells.Item($insertAt, 6).Value = 0
    $ws.Cells.Item($insertAt, 7).Value = 88
    $ws.Cells.Item($insertAt, 8).Value = "https://show.bilibili.com/platform/detail.html?id=85776"
    $ws.Cells.Item($insertAt, 9).Value = "//i2.hdslb.com/bfs/openplatform/202405/Xl4NBnky1715847180514.jpeg"

    # --- fix up the four rows that were pushed down one slot -------------
    # row insertAt+1 : 南宁·AB动漫游戏嘉年华
    $ws.Cells.Item($insertAt + 1, 1).Value = $newRowNumber + 1
    $ws.Cells.Item($insertAt + 1, 6).Value = 1190

    # row insertAt+2 : 横州·第二届海棠动漫游戏嘉年华 (only the sequence number changes)
    $ws.Cells.Item($insertAt + 2, 1).Value = $newRowNumber + 2

    # row insertAt+3 : 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）
    $ws.Cells.Item($insertAt + 3, 1).Value = $newRowNumber + 3
    $ws.Cells.Item($insertAt + 3, 6).Value = 1123

    # row insertAt+4 : 南宁·蔚蓝档案only
    $ws.Cells.Item($insertAt + 4, 1).Value = $newRowNumber + 4
    $ws.Cells.Item($insertAt + 4, 6).Value = 91
}

# "展览" sheet: new row lands at row 10 (its last row before the edit was 13)
$wsExpo = $wb.Worksheets.Item("展览")
Update-ExpoSheet $wsExpo 10

# "全部类型" sheet: one extra pre-existing row, so the new row lands at row 11
$wsAll = $wb.Worksheets.Item("全部类型")
Update-ExpoSheet $wsAll 11
